# New datasets + baseline regression
# Replace every occurrence of the variable name "congenital" with
# "misc_long_term" across all worksheets in the workbook.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Text -eq "congenital") {
            $cell.Value = "misc_long_term"
        }
    }
}
